$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: Monero -> CoreDAO full row change
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.53"
$ws.Range("E51").Value = "  -2.37%  "

# Price / Volume updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.358.61"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.249.64"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.93"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.91"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.246.70"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  +3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.495"
$ws.Range("E12").Value = "  -4.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.93"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.783.65"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.482.31"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.43"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.249.68"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "503.98"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.41"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +3.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.07"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.66"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.14"
$ws.Range("E25").Value = "  +3.12%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.17"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.37"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  +44.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.97"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -4.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.88"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.48"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  +18.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "494.27"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0780"
$ws.Range("E40").Value = "  +14.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0421"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.83"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.980.15"
$ws.Range("E46").Value = "  +5.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.80"
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("E48").Value = "  +5.77%  "
$ws.Range("E49").Value = "  +2.35%  "
